$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 33249.953
$ws.Cells.Item(17, 10).Value = 33249.953
$ws.Cells.Item(17, 12).Value = 99749.859
$ws.Cells.Item(17, 14).Value = -100085.859

# ALC row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 1436.8948
$ws.Cells.Item(38, 9).Value = 300.0909
$ws.Cells.Item(38, 11).Value = 900.2727
$ws.Cells.Item(38, 13).Value = -528.2727

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 2458.5715
$ws.Cells.Item(100, 9).Value = 1937.826
$ws.Cells.Item(100, 10).Value = 4854
$ws.Cells.Item(100, 11).Value = 1937.826
$ws.Cells.Item(100, 12).Value = 4854
$ws.Cells.Item(100, 13).Value = -1396.826
$ws.Cells.Item(100, 14).Value = -5936

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 3451.111
$ws.Cells.Item(113, 10).Value = 3241.4285
$ws.Cells.Item(113, 12).Value = 3241.4285
$ws.Cells.Item(113, 14).Value = -9749.4285

# ALC row 133
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(133, 8).Value = 30000
$ws.Cells.Item(133, 10).Value = 30000
$ws.Cells.Item(133, 12).Value = 30000
$ws.Cells.Item(133, 14).Value = -40120

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 801.94116
$ws.Cells.Item(135, 9).Value = 611.9761999999999
$ws.Cells.Item(135, 10).Value = 1688.4445
$ws.Cells.Item(135, 11).Value = 5507.7858
$ws.Cells.Item(135, 12).Value = 15196.0005
$ws.Cells.Item(135, 13).Value = -2972.7858
$ws.Cells.Item(135, 14).Value = -20266.0005

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 4072.6667
$ws.Cells.Item(137, 9).Value = 4663.4443
$ws.Cells.Item(137, 10).Value = 2300.3333
$ws.Cells.Item(137, 11).Value = 13990.3329
$ws.Cells.Item(137, 12).Value = 6900.999899999999
$ws.Cells.Item(137, 13).Value = -11440.3329
$ws.Cells.Item(137, 14).Value = -12000.9999

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 4435.075
$ws.Cells.Item(138, 9).Value = 2511.75
$ws.Cells.Item(138, 10).Value = 5717.2915
$ws.Cells.Item(138, 11).Value = 7535.25
$ws.Cells.Item(138, 12).Value = 17151.8745
$ws.Cells.Item(138, 13).Value = -2395.25
$ws.Cells.Item(138, 14).Value = -27431.8745

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 410151.03
$ws.Cells.Item(141, 9).Value = 3397.6042
$ws.Cells.Item(141, 10).Value = 1711761.9
$ws.Cells.Item(141, 11).Value = 10192.8126
$ws.Cells.Item(141, 12).Value = 5135285.699999999
$ws.Cells.Item(141, 13).Value = -5012.812600000001
$ws.Cells.Item(141, 14).Value = -5145645.699999999

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5741.31
$ws.Cells.Item(32, 9).Value = 5096.8213
$ws.Cells.Item(32, 10).Value = 9124.875
$ws.Cells.Item(32, 11).Value = 5096.8213
$ws.Cells.Item(32, 12).Value = 9124.875
$ws.Cells.Item(32, 13).Value = -4809.8213
$ws.Cells.Item(32, 14).Value = -9698.875

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 2652.5854
$ws.Cells.Item(132, 9).Value = 1927.3334
$ws.Cells.Item(132, 11).Value = 5782.0002
$ws.Cells.Item(132, 13).Value = -3252.0002

# BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 237.5
$ws.Cells.Item(22, 9).Value = 251.66667
$ws.Cells.Item(22, 10).Value = 223.33333
$ws.Cells.Item(22, 11).Value = 251.66667
$ws.Cells.Item(22, 12).Value = 223.33333
$ws.Cells.Item(22, 13).Value = -78.66667000000001
$ws.Cells.Item(22, 14).Value = -569.3333299999999

# BSM row 51
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(51, 8).Value = 32097.143
$ws.Cells.Item(51, 10).Value = 32097.143
$ws.Cells.Item(51, 12).Value = 32097.143
$ws.Cells.Item(51, 14).Value = -33079.143

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 710
$ws.Cells.Item(94, 9).Value = 637.5
$ws.Cells.Item(94, 10).Value = 1000
$ws.Cells.Item(94, 11).Value = 637.5
$ws.Cells.Item(94, 12).Value = 1000
$ws.Cells.Item(94, 13).Value = -186.5
$ws.Cells.Item(94, 14).Value = -1902

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2060.5518
$ws.Cells.Item(99, 9).Value = 1594.0952
$ws.Cells.Item(99, 10).Value = 3285
$ws.Cells.Item(99, 11).Value = 1594.0952
$ws.Cells.Item(99, 12).Value = 3285
$ws.Cells.Item(99, 13).Value = -96.09519999999998
$ws.Cells.Item(99, 14).Value = -6281

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 1497
$ws.Cells.Item(105, 9).Value = 1374.4445
$ws.Cells.Item(105, 11).Value = 1374.4445
$ws.Cells.Item(105, 13).Value = 372.5554999999999

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2207.8276
$ws.Cells.Item(134, 9).Value = 1954.2439
$ws.Cells.Item(134, 10).Value = 2819.4119
$ws.Cells.Item(134, 11).Value = 5862.7317
$ws.Cells.Item(134, 12).Value = 8458.235700000001
$ws.Cells.Item(134, 13).Value = -3327.7317
$ws.Cells.Item(134, 14).Value = -13528.2357

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2646.9075
$ws.Cells.Item(31, 9).Value = 1638.5897
$ws.Cells.Item(31, 10).Value = 5268.533
$ws.Cells.Item(31, 11).Value = 1638.5897
$ws.Cells.Item(31, 12).Value = 5268.533
$ws.Cells.Item(31, 13).Value = -1343.5897
$ws.Cells.Item(31, 14).Value = -5858.533

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2646.9075
$ws.Cells.Item(34, 9).Value = 1638.5897
$ws.Cells.Item(34, 10).Value = 5268.533
$ws.Cells.Item(34, 11).Value = 1638.5897
$ws.Cells.Item(34, 12).Value = 5268.533
$ws.Cells.Item(34, 13).Value = -1436.5897
$ws.Cells.Item(34, 14).Value = -5672.533

# CRP row 52
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(52, 8).Value = 44000
$ws.Cells.Item(52, 10).Value = 44000
$ws.Cells.Item(52, 12).Value = 44000
$ws.Cells.Item(52, 14).Value = -44588

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 1783.1187
$ws.Cells.Item(132, 9).Value = 1456.1464
$ws.Cells.Item(132, 10).Value = 2527.889
$ws.Cells.Item(132, 11).Value = 4368.439200000001
$ws.Cells.Item(132, 12).Value = 7583.667
$ws.Cells.Item(132, 13).Value = -1838.439200000001
$ws.Cells.Item(132, 14).Value = -12643.667

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 1880.683
$ws.Cells.Item(134, 9).Value = 1216.3636
$ws.Cells.Item(134, 10).Value = 4621
$ws.Cells.Item(134, 11).Value = 3649.0908
$ws.Cells.Item(134, 12).Value = 13863
$ws.Cells.Item(134, 13).Value = -1114.0908
$ws.Cells.Item(134, 14).Value = -18933

# CUL row 37
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 42600
$ws.Cells.Item(37, 10).Value = 42600
$ws.Cells.Item(37, 12).Value = 127800
$ws.Cells.Item(37, 14).Value = -128024

# CUL row 106
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(106, 8).Value = 4000
$ws.Cells.Item(106, 10).Value = 4000
$ws.Cells.Item(106, 12).Value = 12000
$ws.Cells.Item(106, 14).Value = -13892

# CUL row 125
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(125, 8).Value = 1807.5
$ws.Cells.Item(125, 9).Value = 1807.5
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 11).Value = 5422.5
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 13).Value = -502.5
$ws.Cells.Item(125, 14).ClearContents()

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 29540.945
$ws.Cells.Item(102, 9).Value = 1935.8422
$ws.Cells.Item(102, 10).Value = 58679.668
$ws.Cells.Item(102, 11).Value = 1935.8422
$ws.Cells.Item(102, 12).Value = 58679.668
$ws.Cells.Item(102, 13).Value = -313.8422
$ws.Cells.Item(102, 14).Value = -61923.668

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3093.9443
$ws.Cells.Item(132, 9).Value = 2874.6584
$ws.Cells.Item(132, 10).Value = 3785.5386
$ws.Cells.Item(132, 11).Value = 8623.975199999999
$ws.Cells.Item(132, 12).Value = 11356.6158
$ws.Cells.Item(132, 13).Value = -6093.975199999999
$ws.Cells.Item(132, 14).Value = -16416.6158

# GSM row 136
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(136, 8).Value = 17568.375
$ws.Cells.Item(136, 10).Value = 17568.375
$ws.Cells.Item(136, 12).Value = 52705.125
$ws.Cells.Item(136, 14).Value = -57805.125

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 71429384
$ws.Cells.Item(16, 9).Value = 90910056
$ws.Cells.Item(16, 10).Value = 233.33333
$ws.Cells.Item(16, 11).Value = 90910056
$ws.Cells.Item(16, 12).Value = 233.33333
$ws.Cells.Item(16, 13).Value = -90909886
$ws.Cells.Item(16, 14).Value = -573.3333299999999

# LTW row 95
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(95, 8).Value = 20000
$ws.Cells.Item(95, 10).Value = 20000
$ws.Cells.Item(95, 12).Value = 20000
$ws.Cells.Item(95, 14).Value = -25492

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 2351.6047
$ws.Cells.Item(132, 9).Value = 1575.6562
$ws.Cells.Item(132, 10).Value = 4608.909
$ws.Cells.Item(132, 11).Value = 4726.9686
$ws.Cells.Item(132, 12).Value = 13826.727
$ws.Cells.Item(132, 13).Value = -2196.9686
$ws.Cells.Item(132, 14).Value = -18886.727

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 1553.9454
$ws.Cells.Item(136, 9).Value = 1074.9556
$ws.Cells.Item(136, 10).Value = 3709.4
$ws.Cells.Item(136, 11).Value = 3224.8668
$ws.Cells.Item(136, 12).Value = 11128.2
$ws.Cells.Item(136, 13).Value = -674.8667999999998
$ws.Cells.Item(136, 14).Value = -16228.2

# WVR row 97
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(97, 8).Value = 17500
$ws.Cells.Item(97, 10).Value = 17500
$ws.Cells.Item(97, 12).Value = 17500
$ws.Cells.Item(97, 14).Value = -19482

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 7506.816
$ws.Cells.Item(132, 9).Value = 1591.5333
$ws.Cells.Item(132, 11).Value = 4774.5999
$ws.Cells.Item(132, 13).Value = -2244.5999

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 683.87805
$ws.Cells.Item(136, 9).Value = 506.68967
$ws.Cells.Item(136, 10).Value = 1112.0834
$ws.Cells.Item(136, 11).Value = 1520.06901
$ws.Cells.Item(136, 12).Value = 3336.2502
$ws.Cells.Item(136, 13).Value = 1029.93099
$ws.Cells.Item(136, 14).Value = -8436.2502
